$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newStatQuery = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (diag:diagnosis)-->(c)`nOPTIONAL MATCH (f:file)-[*]->(c)`nOPTIONAL MATCH (sf:file)-->(s)`nWITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p`nWHERE demo.breed IN ['Cocker Spaniel']`nRETURN  `n    count(distinct p) AS Programs,`n    count(distinct s) AS Studies,`n    count(distinct c) AS Cases,`n    count(distinct samp) AS Samples,`n    count(distinct f) AS ``Case Files``,`n    count(distinct sf) AS ``Study Files``"

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

$ws.Range("B4").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85
$excel.ActiveWindow.ScrollRow = 3
